$wb = $excel.ActiveWorkbook

# --- Fix Aircraft_scheduling sheet: the aircraft departure date (column D)
# --- was earlier than the landing date (column B) for several rows, making
# --- aircraft appear to depart before they arrived. Push the departure date
# --- to the next day (matches the other rows that already roll over).
$wsAircraft = $wb.Worksheets.Item("Aircraft_scheduling")

$wsAircraft.Range("D2:D6").Value = 45749
$wsAircraft.Range("D9:D10").Value = 45749

# --- Restore / update the selection + active-sheet state so it matches
# --- what was left behind after the manual fix was made in Excel.
$wsWorkPackages = $wb.Worksheets.Item("Work_packages")
$wsStaff = $wb.Worksheets.Item("Staff")

$wsWorkPackages.Activate() | Out-Null
$wsWorkPackages.Range("E8").Select() | Out-Null

$wsStaff.Activate() | Out-Null
$wsStaff.Range("A38").Select() | Out-Null

$wsAircraft.Activate() | Out-Null
$wsAircraft.Range("D11").Select() | Out-Null
